$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.223.07"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "3.551.63"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.33"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.52"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D7").Value = "3.554.00"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.05"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "4.139.61"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000209"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.52"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "3.545.30"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "65.972.02"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.95"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.05"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.20"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.603"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.56"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "3.680.23"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.38"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").Value = "  -6.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.159"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.41"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "3.530.17"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "169.48"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0863"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.895"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -9.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.31"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("E47").Value = "  -7.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.13"
$ws.Range("E48").Value = "  -7.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.20"
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.953"
$ws.Range("E51").Value = "  -4.02%  "
